# Atualizei dados da bibi
# Applies updates to vendas_atipicas: corrects the quantidade_atipica (G)
# values for the existing atypical-sale rows and appends a new atypical
# sale record (row 8) for 2025-06-16.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Corrected quantidade_atipica values on existing rows ---
$ws.Range("G3").Value = -126
$ws.Range("G4").Value = -445
$ws.Range("G5").Value = -65
$ws.Range("G6").Value = -1
$ws.Range("G7").Value = -65

# --- New atypical sale row ---
# A8 and D8 look like a date and a number respectively, but the source
# data stores them as plain text (same as the other rows) - prefix with
# an apostrophe to force text entry, then drop back to the Normal style
# so no stray number-format/quote-prefix style sticks to the cell.
$ws.Range("A8").Value = "'2025-06-16"
$ws.Range("A8").Style = "Normal"

$ws.Range("B8").Value = 2
$ws.Range("C8").Value = "BEMOL S/A"

$ws.Range("D8").Value = "'362396"
$ws.Range("D8").Style = "Normal"

$ws.Range("E8").Value = 13079
$ws.Range("F8").Value = "FONE BLUETOOTH BASIKE TWS FON6694"
$ws.Range("G8").Value = 1
$ws.Range("H8").Value = 1.08
$ws.Range("I8").Value = 0.28
